# further cleaning to metadata
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the index2Sequence value (E7760 -> E7420) for all data rows (2-13)
$ws.Range("K2:K13").Value = "E7420"

# Update the selected range shown in the sheet view from L2:L13 to K2:K13
$ws.Range("K2:K13").Select()

# Convert the boolean roboticLibraryPrep cells (L2:L13) into formula cells
# that evaluate FALSE() while keeping their displayed/cached value as FALSE
for ($r = 2; $r -le 13; $r++) {
    $ws.Cells.Item($r, 12).Formula = "=FALSE()"
}
